# LeetCode Practice Tracker - add "Stack" topic problems
# (Valid Parentheses, Min Stack, Simplify Path ... per commit message)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert 16 blank rows before row 110 so the existing
#    rows 110-112 (String/Valid Anagram, Stack/Valid Parentheses,
#    Queue/Implement Queue using Stacks) shift down to 126-128.
# ---------------------------------------------------------------------
$ws.Rows("110:125").Insert()

# ---------------------------------------------------------------------
# 2. Row 109 was only "108" in column A - fill in the rest of the row
#    (Stack / Valid Parentheses / Easy / Done / date / O(n) / O(n)).
# ---------------------------------------------------------------------
$ws.Range("B109").Value = "Stack"
$ws.Range("C109").Value = "Valid Parentheses"
$ws.Range("D109").Value = "Easy"
$ws.Range("E109").Value = "Done"
$ws.Range("F2").Copy()
$ws.Range("F109").PasteSpecial(-4122)
$ws.Range("F109").Value = 45662
$ws.Range("G109").Value = "O(n)"
$ws.Range("H109").Value = "O(n)"

# ---------------------------------------------------------------------
# 3. New rows 110-123: Stack-topic problems.
# ---------------------------------------------------------------------

# Row 110 - Min Stack (full row, like row 109/123)
$ws.Range("A110").Value = 109
$ws.Range("B110").Value = "Stack"
$ws.Range("C110").Value = "Min Stack"
$ws.Range("D110").Value = "Medium"
$ws.Range("E110").Value = "Done"
$ws.Range("F2").Copy()
$ws.Range("F110").PasteSpecial(-4122)
$ws.Range("F110").Value = 45662
$ws.Range("G110").Value = "O(1)"
$ws.Range("H110").Value = "O(n)"

# Rows 111-122 - just ID / Topic / Problem Name
$problems = @(
    @(111, 110, "Implement Stack using Array / Linked List"),
    @(112, 111, "Implement Stack using Queues"),
    @(113, 112, "Next Greater Element I"),
    @(114, 113, "Next Greater Element II (Circular)"),
    @(115, 114, "Daily Temperatures"),
    @(116, 115, "Evaluate Reverse Polish Notation"),
    @(117, 116, "Largest Rectangle in Histogram"),
    @(118, 117, "Trapping Rain Water (Stack approach)"),
    @(119, 118, "Remove K Digits"),
    @(120, 119, "Decode String"),
    @(121, 120, "Asteroid Collision"),
    @(122, 121, "Online Stock Span")
)
foreach ($p in $problems) {
    $rowNum = $p[0]
    $idVal = $p[1]
    $name = $p[2]
    $ws.Range("A$rowNum").Value = $idVal
    $ws.Range("B$rowNum").Value = "Stack"
    $ws.Range("C$rowNum").Value = $name
}

# Row 123 - Simplify Path (full row; author typed ID 123, skipping 122)
$ws.Range("A123").Value = 123
$ws.Range("B123").Value = "Stack"
$ws.Range("C123").Value = "Simplify Path"
$ws.Range("D123").Value = "Medium"
$ws.Range("E123").Value = "Done"
$ws.Range("F2").Copy()
$ws.Range("F123").PasteSpecial(-4122)
$ws.Range("F123").Value = 45662
$ws.Range("G123").Value = "O(n)"
$ws.Range("H123").Value = "O(n)"

# ---------------------------------------------------------------------
# 4. Rows 124-125: blank placeholder rows, only the ID column filled.
# ---------------------------------------------------------------------
$ws.Range("A124").Value = 124
$ws.Range("A125").Value = 125

# ---------------------------------------------------------------------
# 5. Rows 126-128 already hold the old rows 110-112 content after the
#    insert above - just give them their sequential ID numbers.
# ---------------------------------------------------------------------
$ws.Range("A126").Value = 126
$ws.Range("A127").Value = 127
$ws.Range("A128").Value = 128

# ---------------------------------------------------------------------
# 6. Trailing blank rows 129-130.
# ---------------------------------------------------------------------
$ws.Range("A129").Value = 129
$ws.Range("A130").Value = 130

# ---------------------------------------------------------------------
# 7. Reflect the selected cell after entering the new data (cursor
#    landed on I123 after tabbing through the last new row).
# ---------------------------------------------------------------------
$ws.Range("I123").Select()
